# Add remote firmware update functionality to the register map.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- status / readback registers (rows 110-112, column C) ---
$ws.Range("C110").Value = "REMOTE UPDATE RD ONLY-> status"
$ws.Range("C111").Value = "REMOTE UPDATE RD ONLY-> data low"
$ws.Range("C112").Value = "REMOTE UPDATE RD ONLY-> data high"

# --- remote update control registers (rows 117-126, columns C & D) ---
$ws.Range("C117").Value = "REMOTE UPDATE(0)-> enable/reset"
$ws.Range("D117").Value = "LSB=0 : REMOTE update blocks in reset, commands ignored. Set to 1 to run."

$ws.Range("C118").Value = "REMOTE UPDATE(1) -> write enable to EPCQ fifo "
$ws.Range("D118").Value = "bit 0 write enable, bit 1 toggle fifo write clock"

$ws.Range("C119").Value = "REMOTE UPDATE(2) -> write data to EPCQ fifo (lower 16 bits)"
$ws.Range("D119").Value = "bits 15 to 0"

$ws.Range("C120").Value = "REMOTE UPDATE(3) -> write data to EPCQ fifo (upper 16 bits)"
$ws.Range("D120").Value = "bits 15 to 0"

$ws.Range("C121").Value = "REMOTE UPDATE(4) -> EPCQ command / write mode / clear"
$ws.Range("D121").Value = "cmd - bits 2 to 0 ; mode - bit 8 ; clear - bit 16"

$ws.Range("C122").Value = "REMOTE UPDATE(5) -> write data to EPCQ cmd addr (lower 16 bits)"
$ws.Range("D122").Value = "bits 15 to 0"

$ws.Range("C123").Value = "REMOTE UPDATE(6) -> write data to EPCQ cmd addr (upper 16 bits)"
$ws.Range("D123").Value = "bits 15 to 0"

$ws.Range("C124").Value = "REMOTE UPDATE(7) -> remote update param, toggle_write, reconfig"
$ws.Range("D124").Value = "bits 2 to 0 - param ; bit 8 - toggle_write ; bit 16 - reconfig"

$ws.Range("C125").Value = "REMOTE UPDATE(8) -> remote update data (lower 16 bits)"
$ws.Range("D125").Value = "bits 15 to 0"

$ws.Range("C126").Value = "REMOTE UPDATE(9) -> remote update data (upper 16 bits)"
$ws.Range("D126").Value = "bits 15 to 0"

# --- refresh the view to where the edits were made ---
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("C110").Select()
